$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.087.93'
$ws.Range('E2').Value = '  -0.40%  '

# Row 3
$ws.Range('D3').Value = '1.652.28'
$ws.Range('E3').Value = '  -0.42%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.22%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '217.46'
$ws.Range('E5').Value = '  +0.18%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5277'
$ws.Range('E6').Value = '  +1.98%  '

# Row 7
$ws.Range('E7').Value = '  -0.17%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2598'
$ws.Range('E8').Value = '  -1.58%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06326'
$ws.Range('E9').Value = '  +0.82%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.38'
$ws.Range('E10').Value = '  -2.07%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07788'
$ws.Range('E11').Value = '  +0.23%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.512'
$ws.Range('E12').Value = '  +0.72%  '

# Row 13
$ws.Range('D13').Value = '1.647.11'
$ws.Range('E13').Value = '  -0.70%  '

# Row 14
$ws.Range('D14').Value = '1.879.31'
$ws.Range('E14').Value = '  -0.37%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5495'
$ws.Range('E15').Value = '  +0.50%  '

# Row 16
$ws.Range('D16').Value = '0.0₅8205'
$ws.Range('E16').Value = '  +0.92%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.59'
$ws.Range('E17').Value = '  +0.94%  '

# Row 18
$ws.Range('D18').Value = '26.090.10'
$ws.Range('E18').Value = '  -0.41%  '

# Row 19
$ws.Range('E19').Value = '  -0.21%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.581'
$ws.Range('E20').Value = '  -0.69%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '190.89'
$ws.Range('E21').Value = '  -0.69%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.07'
$ws.Range('E22').Value = '  -0.18%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.036'
$ws.Range('E23').Value = '  +0.40%  '

# Row 24
$ws.Range('E24').Value = '  -0.24%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.14'
$ws.Range('E25').Value = '  +3.39%  '

# Row 26
$ws.Range('E26').Value = '  +1.36%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.234'
$ws.Range('E27').Value = '  -0.66%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '16.05'
$ws.Range('E28').Value = '  -0.66%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.432'
$ws.Range('E29').Value = '  -0.57%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05843'
$ws.Range('E30').Value = '  -1.63%  '

# Row 31
$ws.Range('E31').Value = '  -0.10%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.547'
$ws.Range('E32').Value = '  -0.05%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.267'
$ws.Range('E33').Value = '  -0.55%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.584'
$ws.Range('E34').Value = '  +0.15%  '

# Row 35
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.781'
$ws.Range('E35').Value = '  +0.36%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9450'
$ws.Range('E36').Value = '  -1.60%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.410'
$ws.Range('E37').Value = '  -0.39%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5738'
$ws.Range('E38').Value = '  +0.90%  '

# Row 39
$ws.Range('E39').Value = '  +1.12%  '

# Row 40
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '104.39'
$ws.Range('E40').Value = '  +3.16%  '

# Row 41
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8432'
$ws.Range('E41').Value = '  -1.10%  '

# Row 42
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  -0.15%  '

# Row 43
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.735'
$ws.Range('E43').Value = '  -5.09%  '

# Row 44
$ws.Range('D44').Value = '1.030.14'
$ws.Range('E44').Value = '  +1.73%  '

# Row 45
$ws.Range('D45').Value = '1.794.97'
$ws.Range('E45').Value = '  -0.31%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '57.12'
$ws.Range('E46').Value = '  +1.08%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.002'
$ws.Range('E47').Value = '  -0.57%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4331'
$ws.Range('E48').Value = '  +2.36%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05143'
$ws.Range('E49').Value = '  -0.48%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.798'
$ws.Range('E50').Value = '  -3.27%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.457'
$ws.Range('E51').Value = '  +0.54%  '
